$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above row 648, pushing the existing rows
# 648:693 down to 650:695 (shift down).
$ws.Rows("648:649").Insert(-4121)

# Row 648 (new) - Primera, Provincia de Santiago
$ws.Cells.Item(648, 1).Value = 3
$ws.Cells.Item(648, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(648, 3).Value = "Coquimbo"
$ws.Cells.Item(648, 4).Value = 44714
$ws.Cells.Item(648, 5).Value = 5
$ws.Cells.Item(648, 6).Value = 100112023
$ws.Cells.Item(648, 7).Value = "Brócoli"
$ws.Cells.Item(648, 8).Value = "Sin especificar"
$ws.Cells.Item(648, 9).Value = "Primera"
$ws.Cells.Item(648, 10).Value = 2800
$ws.Cells.Item(648, 11).Value = 850
$ws.Cells.Item(648, 12).Value = 900
$ws.Cells.Item(648, 13).Value = 871
$ws.Cells.Item(648, 14).Value = "`$/unidad"
$ws.Cells.Item(648, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(648, 16).Value = 871
$ws.Cells.Item(648, 17).Value = 1
$ws.Cells.Item(648, 18).Value = "Hortaliza"

# Row 649 (new) - Segunda, Provincia de Santiago
$ws.Cells.Item(649, 1).Value = 3
$ws.Cells.Item(649, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(649, 3).Value = "Coquimbo"
$ws.Cells.Item(649, 4).Value = 44714
$ws.Cells.Item(649, 5).Value = 5
$ws.Cells.Item(649, 6).Value = 100112023
$ws.Cells.Item(649, 7).Value = "Brócoli"
$ws.Cells.Item(649, 8).Value = "Sin especificar"
$ws.Cells.Item(649, 9).Value = "Segunda"
$ws.Cells.Item(649, 10).Value = 1100
$ws.Cells.Item(649, 11).Value = 600
$ws.Cells.Item(649, 12).Value = 600
$ws.Cells.Item(649, 13).Value = 600
$ws.Cells.Item(649, 14).Value = "`$/unidad"
$ws.Cells.Item(649, 15).Value = "Provincia de Santiago"
$ws.Cells.Item(649, 16).Value = 600
$ws.Cells.Item(649, 17).Value = 1
$ws.Cells.Item(649, 18).Value = "Hortaliza"
